# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# Rushing sheet updates (row 2 = B.Mayfield, row 3 = C.Keenum, row 11 = J.Landry)
$rushing.Range("C2").Value = 4
$rushing.Range("E2").Value = 5

$rushing.Range("E3").Value = 3
$rushing.Range("F3").Value = 2

$rushing.Range("C11").Value = 1
$rushing.Range("D11").Value = 2
$rushing.Range("E11").Value = 2
$rushing.Range("F11").Value = 4

# Receiving sheet updates (row 8 = J.Landry, row 13 = A.Hooper)
$receiving.Range("C8").Value = 58
$receiving.Range("D8").Value = 47
$receiving.Range("E8").Value = 22
$receiving.Range("F8").Value = 11
$receiving.Range("G8").Value = 7
$receiving.Range("H8").Value = 3

$receiving.Range("C13").Value = 48
$receiving.Range("D13").Value = 30
$receiving.Range("E13").Value = 5
$receiving.Range("F13").Value = 3
$receiving.Range("G13").Value = 10
$receiving.Range("H13").Value = 6
